$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.290.68"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "3.524.39"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.45%  "
$ws.Range("E6").Value = "  -6.08%  "
$ws.Range("E7").Value = "  +3.88%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E10").Value = "  +7.36%  "
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("E12").Value = "  +4.01%  "
$ws.Range("D14").Value = "4.081.41"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "3.526.04"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "66.196.97"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "417.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.83%  "
$ws.Range("E22").Value = "  +10.41%  "
$ws.Range("E23").Value = "  +5.52%  "
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("E25").Value = "  +11.76%  "
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  -1.77%  "
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "608.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E36").Value = "  +8.60%  "
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +9.69%  "
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").Value = "3.256.31"
$ws.Range("E42").Value = "  +7.63%  "
$ws.Range("E43").Value = "  +4.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  -6.77%  "
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.09%  "
